$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.162.11"
$ws.Range("E2").Value = "  +7.07%  "

# Row 3
$ws.Range("D3").Value = "3.599.76"
$ws.Range("E3").Value = "  +3.83%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "415.12"
$ws.Range("E5").Value = "  +0.24%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.50"
$ws.Range("E6").Value = "  -0.61%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.651"
$ws.Range("E7").Value = "  +3.89%  "

# Row 8
$ws.Range("D8").Value = "3.593.14"
$ws.Range("E8").Value = "  +3.76%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.780"
$ws.Range("E10").Value = "  +7.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.175"
$ws.Range("E11").Value = "  +18.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000337"
$ws.Range("E12").Value = "  +54.38%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.39"
$ws.Range("E13").Value = "  -0.36%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.88"
$ws.Range("E14").Value = "  +2.81%  "

# Row 15
$ws.Range("D15").Value = "4.166.02"
$ws.Range("E15").Value = "  +3.66%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.28"
$ws.Range("E17").Value = "  -1.08%  "

# Row 18
$ws.Range("D18").Value = "3.616.78"
$ws.Range("E18").Value = "  +3.17%  "

# Row 19
$ws.Range("E19").Value = "  +5.37%  "

# Row 20
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "67.031.84"
$ws.Range("E20").Value = "  +6.87%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.30"
$ws.Range("E21").Value = "  -2.94%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "448.45"
$ws.Range("E22").Value = "  -2.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "89.19"
$ws.Range("E23").Value = "  -1.51%  "

# Row 24
$ws.Range("E24").Value = "  -3.33%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.10"
$ws.Range("E25").Value = "  -1.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.34"
$ws.Range("E26").Value = "  +0.76%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.93"
$ws.Range("E27").Value = "  -7.52%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "35.21"
$ws.Range("E28").Value = "  +5.70%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.88"
$ws.Range("E29").Value = "  +1.80%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.38"
$ws.Range("E30").Value = "  +3.34%  "

# Row 31
$ws.Range("E31").Value = "  +3.58%  "

# Row 32
$ws.Range("E32").Value = "  +4.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.35"
$ws.Range("E33").Value = "  -2.81%  "

# Row 34
$ws.Range("E34").Value = "  -3.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "40.07"
$ws.Range("E35").Value = "  -1.59%  "

# Row 36
$ws.Range("E36").Value = "  +0.07%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.76"
$ws.Range("E37").Value = "  -2.90%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0493"
$ws.Range("E38").Value = "  +0.39%  "

# Row 39
$ws.Range("D39").Value = "0.0₃0725"
$ws.Range("E39").Value = "  +30.11%  "

# Row 40
$ws.Range("E40").Value = "  +9.91%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E41").Value = "  -0.17%  "

# Row 42
$ws.Range("E42").Value = "  -2.10%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "149.47"
$ws.Range("E43").Value = "  +1.24%  "

# Row 44
$ws.Range("E44").Value = "  +1.98%  "

# Row 45
$ws.Range("B45").Value = "LidoDAOToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  -1.78%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.314"
$ws.Range("E46").Value = "  -1.99%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.30"
$ws.Range("E47").Value = "  -1.10%  "

# Row 48
$ws.Range("E48").Value = "  -4.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.29"
$ws.Range("E49").Value = "  -4.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "15.67"
$ws.Range("E50").Value = "  -4.37%  "

# Row 51
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "114.46"
$ws.Range("E51").Value = "  +5.46%  "
